$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.985.87"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.601.93"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.60"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.484"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.13"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.57"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.604.98"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.986.94"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.24"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.87"
$ws.Range("E20").Value = "  +8.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.23"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.81"
$ws.Range("E24").Value = "  +8.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.37"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.121"
$ws.Range("E27").Value = "  -6.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.15"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.11"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0165"
$ws.Range("E36").Value = "  +10.94%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.127.71"
$ws.Range("E37").Value = "  +4.44%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.793"
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.491"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.737.55"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.77"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.46"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.407"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0924"
$ws.Range("E51").Value = "  -17.53%  "
